$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MISSION.TYPE")

# The "Code" column (A) for rows 10-47 contained values prefixed with a
# leading "/" (e.g. "/SAV/ASC"). Strip that leading slash from each one.
$ws.Range("A10").Value2 = "SAV/ASC"
$ws.Range("A11").Value2 = "FSTT/TA/FR_DIFFICULTACCSS"
$ws.Range("A12").Value2 = "GEN/RECVRY"
$ws.Range("A13").Value2 = "RSC/SAR/FR_VLN"
$ws.Range("A14").Value2 = "SAV/ASC/FR_PPL/LIFT"
$ws.Range("A15").Value2 = "SAV/RHD"
$ws.Range("A16").Value2 = "SAV/RTA"
$ws.Range("A17").Value2 = "SAV/SARCSL"
$ws.Range("A18").Value2 = "FFST/FR_FIRE"
$ws.Range("A19").Value2 = "FSTT/RRHAZ/FR_CO"
$ws.Range("A20").Value2 = "FSTT/TA/FR_ANI/DGR"
$ws.Range("A21").Value2 = "FSTT/TA/FR_ANI/INJ"
$ws.Range("A22").Value2 = "INT/RECCE/FR_CBRNHZ"
$ws.Range("A23").Value2 = "INT/RECCE/FR_DIS/LNDSLD"
$ws.Range("A24").Value2 = "INT/RECCE/FR_DIS/SDCLPS"
$ws.Range("A25").Value2 = "INT/RECCE/FR_FLD"
$ws.Range("A26").Value2 = "INT/RECCE/FR_SMLL"
$ws.Range("A27").Value2 = "FSTT/TA/FR_DRG/MIND"
$ws.Range("A28").Value2 = "FR_MED/REGLTN"
$ws.Range("A29").Value2 = "GEN/SUPRTN"
$ws.Range("A30").Value2 = "REC/PRVCNP"
$ws.Range("A31").Value2 = "RSC/MEDEVC"
$ws.Range("A32").Value2 = "SAV/AR/FR_CNT"
$ws.Range("A33").Value2 = "SAV/AR/FR_MED"
$ws.Range("A34").Value2 = "SAV/AR/FR_PARAMD"
$ws.Range("A35").Value2 = "SAV/AR/FR_PPL/GRP"
$ws.Range("A36").Value2 = "SAV/AR/FR_PSYPHY"
$ws.Range("A37").Value2 = "SAV/ASC"
$ws.Range("A38").Value2 = "SAV/ASC"
$ws.Range("A39").Value2 = "INT/RECCE"
$ws.Range("A40").Value2 = "SAV"
$ws.Range("A41").Value2 = "FSTT/TA/FR_TRNSP/AMB"
$ws.Range("A42").Value2 = "GEN/TRNSP"
$ws.Range("A43").Value2 = "GEN/TRNSP/FR_SECNDRY"
$ws.Range("A44").Value2 = "GEN/TRNSPN"
$ws.Range("A45").Value2 = "OPR/LOG"
$ws.Range("A46").Value2 = "SAV/AR/FR_PPL/OBS"
$ws.Range("A47").Value2 = "FSTT/TA/FR_CLRACCSS"

# Update the active selection on the sheet to reflect where the editor's
# cursor ended up after making the change.
$ws.Activate()
$ws.Range("H15").Select()
